$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values are written with a leading apostrophe so that
# numeric-looking strings (e.g. "561.95") stay text instead of being
# auto-converted to numbers (which would drop formatting such as a
# trailing zero in "87.20").

$ws.Range("D2").Value = "'65.420.30"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "'3.408.01"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'561.95"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").Value = "'176.16"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").Value = "'3.398.94"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +12.62%  "
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("D12").Value = "'54.99"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("E13").Value = "  +5.57%  "
$ws.Range("D14").Value = "'9.16"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "'3.949.83"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "'18.39"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "'3.399.74"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "'65.454.57"
$ws.Range("D20").Value = "'11.92"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'473.88"
$ws.Range("E22").Value = "  +15.25%  "
$ws.Range("D23").Value = "'5.07"
$ws.Range("E23").Value = "  +16.69%  "
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").Value = "'87.20"
$ws.Range("E25").Value = "  +5.09%  "
$ws.Range("D26").Value = "'13.49"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("E28").Value = "  +6.46%  "
$ws.Range("D29").Value = "'8.88"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("D30").Value = "'31.29"
$ws.Range("E30").Value = "  +7.64%  "
$ws.Range("D31").Value = "'6.77"
$ws.Range("E31").Value = "  +6.34%  "
$ws.Range("D33").Value = "'62.56"
$ws.Range("E33").Value = "  +7.95%  "
$ws.Range("D34").Value = "'575.85"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "'35.93"
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("D40").Value = "'0.0₃0761"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").Value = "'3.092.98"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  +4.72%  "
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("E47").Value = "  +5.65%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D50").Value = "'137.35"
$ws.Range("E50").Value = "  +3.97%  "
$ws.Range("D51").Value = "'8.34"
$ws.Range("E51").Value = "  +3.31%  "
